# Update the user-story wireframe sheet:
#  - the title cell (B2) gets a new, shorter title
#  - the "En tant que" description (C3) is reworded
# then move the active selection to C3 (matches the saved UI state in the
# target workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Afficher le details d'une tâche"
$ws.Range("C3").Value = "je suis un utilisateur connecté et que je suis sur la page de calendrier"

$ws.Range("C3").Select() | Out-Null
